$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet
$ws.Name = "MV_Testdaten"

# 2) Remove the existing table definition so we can rebuild it cleanly once the
#    new column has been inserted (keeps header-name -> column mapping correct).
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# 3) Insert a new column at E ("Buchungssprache") between "Rolle" and
#    "Zulassungsland", shifting the remaining columns right.
$ws.Range("E1").EntireColumn.Insert()

# Give the new column the same width as its neighbours (C:D).
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# 4) Populate the new column's header + data.
$ws.Range("E1").Value2 = "Buchungssprache"
$ws.Range("E2").Value2 = "Deutsch"
$ws.Range("E3").Value2 = "Polski"
$ws.Range("E4").Value2 = "Deutsch"

# 5) Rebuild the table over the new, wider range.
$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:K9"), [System.Reflection.Missing]::Value, 1)
$newLo.Name = "Table1"
$newLo.TableStyle = "TableStyleMedium9"

# 6) Restore the active selection to E1.
$ws.Range("E1").Select() | Out-Null
